$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.021.51"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").Value = "1.923.83"
$ws.Range("E3").Value = "  +0.89%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.80"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4589"
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3817"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07743"
$ws.Range("E9").Value = "  -0.09%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9791"
$ws.Range("E10").Value = "  -0.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.58"
$ws.Range("E11").Value = "  +2.22%  "
$ws.Range("D12").Value = "1.960.02"
$ws.Range("E12").Value = "  +0.44%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.703"
$ws.Range("E13").Value = "  +0.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.973"
$ws.Range("E14").Value = "  -0.40%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07003"
$ws.Range("E15").Value = "  -0.69%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "84.83"
$ws.Range("E16").Value = "  +0.70%  "
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009497"
$ws.Range("E18").Value = "  -0.82%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.70"
$ws.Range("E19").Value = "  -0.45%  "
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("D21").Value = "29.075.04"
$ws.Range("E21").Value = "  +0.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.355"
$ws.Range("E22").Value = "  +0.26%  "
$ws.Range("E23").Value = "  +0.86%  "
$ws.Range("D24").Value = "2.175.41"
$ws.Range("E24").Value = "  +0.21%  "
$ws.Range("E25").Value = "  -1.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "157.99"
$ws.Range("E26").Value = "  +0.81%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.01"
$ws.Range("E27").Value = "  -0.71%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.623"
$ws.Range("E28").Value = "  +0.22%  "
$ws.Range("E29").Value = "  -0.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.829"
$ws.Range("E30").Value = "  -0.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09323"
$ws.Range("E31").Value = "  +0.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8615"
$ws.Range("E32").Value = "  -0.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.097"
$ws.Range("E33").Value = "  -0.30%  "
$ws.Range("E34").Value = "  -0.71%  "
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05701"
$ws.Range("E36").Value = "  -0.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.148"
$ws.Range("E37").Value = "  +0.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.003"
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02050"
$ws.Range("E39").Value = "  +0.48%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.095"
$ws.Range("E40").Value = "  +13.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.448"
$ws.Range("E41").Value = "  -0.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5510"
$ws.Range("E42").Value = "  -0.58%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1755"
$ws.Range("E43").Value = "  -0.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "9.339"
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.000002808"
$ws.Range("E45").Value = "  +7.58%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.182"
$ws.Range("E46").Value = "  +3.99%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5184"
$ws.Range("E47").Value = "  -0.48%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06942"
$ws.Range("E48").Value = "  +1.84%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "11.20"
$ws.Range("E49").Value = "  -1.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "110.85"
$ws.Range("E50").Value = "  -0.77%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.760"
$ws.Range("E51").Value = "  -1.10%  "
